$d = $word.ActiveDocument

$pairs = @(
    @("94-4=", "29+49="),
    @("35+33=", "42+10="),
    @("56-50=", "35+48="),
    @("78+19=", "17+81="),
    @("53-25=", "60-54="),
    @("6+10=", "71-45="),
    @("55+20=", "76-72="),
    @("8+39=", "73+18="),
    @("14-10=", "10+50="),
    @("17+62=", "89-16="),
    @("77+14=", "51+31="),
    @("11+25=", "97-10="),
    @("23+2=", "79-7="),
    @("41+12=", "39+18="),
    @("90-28=", "28+59="),
    @("99-15=", "11+34="),
    @("37+8=", "89-12="),
    @("24+21=", "11+68="),
    @("45-44=", "80-77="),
    @("36+6=", "95-74="),
    @("30+30=", "61-35="),
    @("20+75=", "73+6="),
    @("8+18=", "51+7="),
    @("41+8=", "62-60="),
    @("65+34=", "26+52="),
    @("16+23=", "51-22="),
    @("60+1=", "48-1="),
    @("17+66=", "20-17="),
    @("46-32=", "33+5="),
    @("76-26=", "88-55="),
    @("38+58=", "47-33="),
    @("44-14=", "28+34="),
    @("3+6=", "25+74="),
    @("90-87=", "42-30="),
    @("88-40=", "80-5="),
    @("67+21=", "70-56="),
    @("31+18=", "15+23="),
    @("74-57=", "14-8="),
    @("18+45=", "56-22="),
    @("72-57=", "98-94="),
    @("82-57=", "69-35="),
    @("41-39=", "54-34="),
    @("69-24=", "68-36="),
    @("91+5=", "47-26="),
    @("31-25=", "65+30="),
    @("43+41=", "66-58="),
    @("94-50=", "96-53="),
    @("53+39=", "69-57="),
    @("44+31=", "78+2="),
    @("91-49=", "51-46="),
    @("21-10=", "63-15="),
    @("61+32=", "61+3="),
    @("72-45=", "84+0="),
    @("60+0=", "36+15="),
    @("92-49=", "94-79="),
    @("27+8=", "51+29="),
    @("72-12=", "65-0="),
    @("91-20=", "39+7="),
    @("0+66=", "47-9="),
    @("78-50=", "55+16="),
    @("44+45=", "58+21="),
    @("26+56=", "60+37="),
    @("63-43=", "38-8="),
    @("76+16=", "92-43="),
    @("9+62=", "90-85="),
    @("12+38=", "72+21="),
    @("80-46=", "0+43="),
    @("96-62=", "32+9="),
    @("65-6=", "79-5="),
    @("11-4=", "97-11="),
    @("38-25=", "80+15="),
    @("35+14=", "89-13="),
    @("35-30=", "90-2="),
    @("37-8=", "89-14="),
    @("4+75=", "86-2="),
    @("71+26=", "22+25="),
    @("75-11=", "55+38="),
    @("70-21=", "49+4="),
    @("1+47=", "15+14="),
    @("16+12=", "44-38="),
    @("46+32=", "23+16="),
    @("54+16=", "77-50="),
    @("9+37=", "96-85="),
    @("32+16=", "94-38="),
    @("22+76=", "2+29="),
    @("74-23=", "37-12="),
    @("97-91=", "83+9="),
    @("7+10=", "72+8="),
    @("83-42=", "0+67="),
    @("81+18=", "89+4="),
    @("12-7=", "90-89="),
    @("61-59=", "53-38="),
    @("71-54=", "87-54="),
    @("54-44=", "56+23="),
    @("44+16=", "55+27="),
    @("73+2=", "27-26="),
    @("13+85=", "60-42="),
    @("25+55=", "48+46="),
    @("17-1=", "7+52="),
    @("75-41=", "1+73="),
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Output "Replaced $($pairs.Count) equations"